$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Milestone_1")
$ws2 = $wb.Worksheets.Item("Milestone_2")

# Select the worksheet that is being edited (Milestone_2 / sheet2)
$ws2.Select()

# Copy the existing "Responsible Team Member" cell format (centered, wrapped,
# bordered style already used on the Milestone_1 sheet) onto the same column
# of rows 12-19 on Milestone_2 so the new style matches without inventing a
# brand-new cell style.
$ws1.Range("B12").Copy()
$ws2.Range("B12:B19").PasteSpecial(-4122)

# Fill in the "Responsible Team Member" column for the milestone 2 task list.
# Values are entered in this particular order so that any brand-new shared
# strings get appended to the shared-strings table in the same order as the
# original edit.
$ws2.Range("B15").Value = "Sodara"
$ws2.Range("B16").Value = "Joe/Shay"
$ws2.Range("B18").Value = "Megan/Joe"
$ws2.Range("B14").Value = "Joe/Megan"
$ws2.Range("B12").Value = "Shay"
$ws2.Range("B13").Value = "Joe"
$ws2.Range("B17").Value = "Megan"
$ws2.Range("B19").Value = "Joe"

# Update the view's selected cell to match where the editor left off.
$ws2.Range("D13").Select()
